$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Re-use the existing date/time style (s="1") from the last data row
# so no new style is introduced for the new A-column timestamps.
$ws.Range("A118").Copy()
$ws.Range("A119:A120").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# ---- Row 119 ----
$ws.Range("A119").Value = 45460.2916666667
$ws.Range("B119").Value = 0
$ws.Range("C119").Value = 2.46000003814697
$ws.Range("D119").Value = 2.46000003814697
$ws.Range("E119").Value = 2.46000003814697
$ws.Range("F119").Value = 2.46000003814697

# adj_close is stored as text (matching the source R-script export).
# A leading apostrophe forces Excel to keep the numeric-looking value as
# text instead of auto-converting it to a number.
$ws.Range("G119").Value = "'2.46000003814697"
$ws.Range("G119").Style = "Normal"

$ws.Range("H119").Value = "LS.MI"

# ---- Row 120 ----
$ws.Range("A120").Value = 45461.6278356481
$ws.Range("B120").Value = 7000
$ws.Range("C120").Value = 2.39000010490417
$ws.Range("D120").Value = 2.32999992370605
$ws.Range("E120").Value = 2.39000010490417
$ws.Range("F120").Value = 2.33999991416931

$ws.Range("G120").Value = "'2.33999991416931"
$ws.Range("G120").Style = "Normal"

$ws.Range("H120").Value = "LS.MI"
